$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 960.3
$ws.Range("J112").Value = 991.3929000000001
$ws.Range("L112").Value = 2974.1787
$ws.Range("N112").Value = -5190.1787

$ws.Range("H137").Value = 1578.3077
$ws.Range("I137").Value = 1563.6666
$ws.Range("K137").Value = 4690.9998
$ws.Range("M137").Value = -2140.9998

$ws.Range("H138").Value = 2250.161
$ws.Range("I138").Value = 1675.3182
$ws.Range("K138").Value = 5025.9546
$ws.Range("M138").Value = 114.0454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 966
$ws.Range("I21").Value = 966
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 966
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -592
$ws.Range("N21").ClearContents()

$ws.Range("H61").Value = 1214.3478
$ws.Range("I61").Value = 1091.9524
$ws.Range("J61").Value = 2499.5
$ws.Range("K61").Value = 1091.9524
$ws.Range("L61").Value = 2499.5
$ws.Range("M61").Value = -879.9523999999999
$ws.Range("N61").Value = -2923.5

$ws.Range("H74").Value = 2566
$ws.Range("I74").Value = 1217.1875
$ws.Range("J74").Value = 5263.625
$ws.Range("K74").Value = 1217.1875
$ws.Range("L74").Value = 5263.625
$ws.Range("M74").Value = -343.1875
$ws.Range("N74").Value = -7011.625

$ws.Range("H77").Value = 2566
$ws.Range("I77").Value = 1217.1875
$ws.Range("J77").Value = 5263.625
$ws.Range("K77").Value = 6085.9375
$ws.Range("L77").Value = 26318.125
$ws.Range("M77").Value = -1717.9375
$ws.Range("N77").Value = -35054.125

$ws.Range("H132").Value = 5562.6665
$ws.Range("I132").Value = 7564.8335
$ws.Range("J132").Value = 3560.5
$ws.Range("K132").Value = 22694.5005
$ws.Range("L132").Value = 10681.5
$ws.Range("M132").Value = -20164.5005
$ws.Range("N132").Value = -15741.5

$ws.Range("H136").Value = 1214.3478
$ws.Range("I136").Value = 1091.9524
$ws.Range("J136").Value = 2499.5
$ws.Range("K136").Value = 3275.857199999999
$ws.Range("L136").Value = 7498.5
$ws.Range("M136").Value = -725.8571999999995
$ws.Range("N136").Value = -12598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2581.3635
$ws.Range("I134").Value = 2707.8076
$ws.Range("K134").Value = 8123.4228
$ws.Range("M134").Value = -5588.4228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17386.953
$ws.Range("I31").Value = 28931.889
$ws.Range("J31").Value = 3055.3103
$ws.Range("K31").Value = 28931.889
$ws.Range("L31").Value = 3055.3103
$ws.Range("M31").Value = -28636.889
$ws.Range("N31").Value = -3645.3103

$ws.Range("H34").Value = 17386.953
$ws.Range("I34").Value = 28931.889
$ws.Range("J34").Value = 3055.3103
$ws.Range("K34").Value = 28931.889
$ws.Range("L34").Value = 3055.3103
$ws.Range("M34").Value = -28729.889
$ws.Range("N34").Value = -3459.3103

$ws.Range("H58").Value = 9390.156000000001
$ws.Range("I58").Value = 1580.1177
$ws.Range("J58").Value = 18241.533
$ws.Range("K58").Value = 1580.1177
$ws.Range("L58").Value = 18241.533
$ws.Range("M58").Value = -1377.1177
$ws.Range("N58").Value = -18647.533

$ws.Range("H92").Value = 26225
$ws.Range("J92").Value = 26225
$ws.Range("L92").Value = 26225
$ws.Range("N92").Value = -31217

$ws.Range("H120").Value = 23645.455
$ws.Range("J120").Value = 23645.455
$ws.Range("L120").Value = 23645.455
$ws.Range("N120").Value = -30903.455

$ws.Range("H132").Value = 36588580
$ws.Range("I132").Value = 34485724
$ws.Range("J132").Value = 41670480
$ws.Range("K132").Value = 103457172
$ws.Range("L132").Value = 125011440
$ws.Range("M132").Value = -103454642
$ws.Range("N132").Value = -125016500

$ws.Range("H134").Value = 1675.6
$ws.Range("I134").Value = 1839.3334
$ws.Range("K134").Value = 5518.0002
$ws.Range("M134").Value = -2983.0002

$ws.Range("H136").Value = 9390.156000000001
$ws.Range("I136").Value = 1580.1177
$ws.Range("J136").Value = 18241.533
$ws.Range("K136").Value = 4740.3531
$ws.Range("L136").Value = 54724.599
$ws.Range("M136").Value = -2190.3531
$ws.Range("N136").Value = -59824.599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 333334240
$ws.Range("J25").Value = 500000600
$ws.Range("L25").Value = 1500001800
$ws.Range("N25").Value = -1500002138

$ws.Range("H30").Value = 333334240
$ws.Range("J30").Value = 500000600
$ws.Range("L30").Value = 1500001800
$ws.Range("N30").Value = -1500002004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 54248.242
$ws.Range("I70").Value = 88389.336
$ws.Range("J70").Value = 6049.0586
$ws.Range("K70").Value = 88389.336
$ws.Range("L70").Value = 6049.0586
$ws.Range("M70").Value = -88119.336
$ws.Range("N70").Value = -6589.0586

$ws.Range("H73").Value = 54248.242
$ws.Range("I73").Value = 88389.336
$ws.Range("J73").Value = 6049.0586
$ws.Range("K73").Value = 88389.336
$ws.Range("L73").Value = 6049.0586
$ws.Range("M73").Value = -87453.336
$ws.Range("N73").Value = -7921.0586

$ws.Range("H132").Value = 2608.3333
$ws.Range("I132").Value = 1862.625
$ws.Range("K132").Value = 5587.875
$ws.Range("M132").Value = -3057.875

$ws.Range("H141").Value = 50359.4
$ws.Range("J141").Value = 50359.4
$ws.Range("L141").Value = 50359.4
$ws.Range("N141").Value = -60719.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 941.17645
$ws.Range("J22").Value = 875.0625
$ws.Range("L22").Value = 875.0625
$ws.Range("N22").Value = -1465.0625

$ws.Range("H27").Value = 941.17645
$ws.Range("J27").Value = 875.0625
$ws.Range("L27").Value = 875.0625
$ws.Range("N27").Value = -1089.0625

$ws.Range("H97").Value = 19122
$ws.Range("J97").Value = 19122
$ws.Range("L97").Value = 19122
$ws.Range("N97").Value = -21104

$ws.Range("H132").Value = 1643.9333
$ws.Range("I132").Value = 1424.1904
$ws.Range("K132").Value = 4272.5712
$ws.Range("M132").Value = -1742.5712

$ws.Range("H136").Value = 1471.6923
$ws.Range("I136").Value = 1450.6086
$ws.Range("K136").Value = 4351.825800000001
$ws.Range("M136").Value = -1801.825800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 251000
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1707

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H80").Value = 3959393.5
$ws.Range("J80").Value = 3959393.5
$ws.Range("L80").Value = 3959393.5
$ws.Range("N80").Value = -3961389.5

$ws.Range("H83").Value = 3959393.5
$ws.Range("J83").Value = 3959393.5
$ws.Range("L83").Value = 11878180.5
$ws.Range("N83").Value = -11888164.5

$ws.Range("H132").Value = 3009.48
$ws.Range("I132").Value = 3039.9473
$ws.Range("J132").Value = 2913
$ws.Range("K132").Value = 9119.841899999999
$ws.Range("L132").Value = 8739
$ws.Range("M132").Value = -6589.841899999999
$ws.Range("N132").Value = -13799

$ws.Range("H136").Value = 1140.3636
$ws.Range("I136").Value = 885.7143
$ws.Range("J136").Value = 1586
$ws.Range("K136").Value = 2657.1429
$ws.Range("L136").Value = 4758
$ws.Range("M136").Value = -107.1428999999998
$ws.Range("N136").Value = -9858
